# Atualização dados 15 e 16/10
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - HUSE
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 0

# Row 3 - Hospital Cirurgia - SUS
$ws.Range("C3").Value = 8

# Row 6 - Hospital Regional - Estancia
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = 1

# Row 7 - Hospital N. Sra da Conceicao
$ws.Range("C7").Value = 2

# Row 11 - Hospital Sao Lucas/Adulto
$ws.Range("C11").Value = 2

# Row 12 - Hospital Sao Lucas/Pediatrico
$ws.Range("C12").Value = 1

# Row 13 - Hospital Gabriel Soares/Adulto
$ws.Range("C13").Value = 0

# Row 17 - Hospital Primavera
$ws.Range("C17").Value = 4

# Update selection to match new active cell
$ws.Range("D8").Select()
